$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update summary figures
$ws.Range("E11").Value = 60000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Remove the second worker's three rows (17-19) entirely
$ws.Range("A17:J19").EntireRow.Delete()
